$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "286.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.31%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.73%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.198"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.66%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06963"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.71%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.420"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.67%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.559"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.74%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.415"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.68%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8967"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.42%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1600"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.25%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07697"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "26.10%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07654"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.99%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02925"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.11%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08996"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.10%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001601"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "3.13%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006516"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.46%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006130"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.47%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.461"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.20%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.43%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3228"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.83%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1330"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.79%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.017"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.63%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.72%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04518"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.12%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001210"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "3.36%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004232"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-5.10%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001169"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.66%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001639"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2.48%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04310"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.79%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006934"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.23%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1242"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.01%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002078"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.78%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01160"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.27%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005847"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.85%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01306"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "1.09%"
